$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 18 first so its shared string ("AVR compiles...") is registered
#     before row 17's ("HSERPRINT...."), matching the author's save order
#     (sharedStrings index 24 = AVR..., index 25 = HSERPRINT...).
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "OPEN"
$ws.Range("D18").Value = "AVR compiles when it should not.  See https://sourceforge.net/p/gcbasic/discussion/596084/thread/e58866dc/#5e0f"

# --- Row 17 (#16 - Closed, HSERPRINT/LONGs fix) ---
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "CLOSED"
$ws.Range("D17").Value = "HSERPRINT not handling LONGs correctly." + [char]10 + "Updated USART.H to handle LONGs correctly."

# Row 17's description cell wraps text like the other description cells in
# column D (style index 2 in the original workbook).
$ws.Range("D17").WrapText = $true
$ws.Range("D17").HorizontalAlignment = -4131
$ws.Range("D17").VerticalAlignment = -4160

# Match the explicit row height recorded in the sheet for row 17.
$ws.Rows.Item(17).RowHeight = 30

# Update the view: scroll so row 13 is at the top and select D18, mirroring
# the workbook's saved cursor/viewport position after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D18").Select()
